$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without Excel auto-converting
# number-looking / percent-looking strings into numeric cells, and without
# mutating the cell's NumberFormat/Style (which a straight ".Value = ..." or
# ".NumberFormat = '@'" assignment would do). We build the literal as a text
# formula (="...") then collapse the formula down to its computed value via
# Copy + PasteSpecial(xlPasteValues), leaving a plain static text cell behind.
function Set-CellText {
    param($sheet, [string]$cellRef, [string]$text)
    $cell = $sheet.Range($cellRef)
    $cell.Formula = ('="' + $text + '"')
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

$excel.CutCopyMode = 0

Set-CellText $ws 'D2' '302.38'
Set-CellText $ws 'E2' '0.67%'
Set-CellText $ws 'D3' '32.23'
Set-CellText $ws 'E3' '1.61%'
Set-CellText $ws 'D4' '4.980'
Set-CellText $ws 'E4' '-3.27%'
Set-CellText $ws 'D5' '0.07907'
Set-CellText $ws 'E5' '-2.78%'
Set-CellText $ws 'D6' '2.136'
Set-CellText $ws 'E6' '-13.98%'
Set-CellText $ws 'D7' '7.834'
Set-CellText $ws 'E7' '0.54%'
Set-CellText $ws 'D8' '3.811'
Set-CellText $ws 'E8' '-1.91%'
Set-CellText $ws 'D9' '0.9249'
Set-CellText $ws 'E9' '0.10%'
Set-CellText $ws 'D10' '0.1752'
Set-CellText $ws 'E10' '-0.31%'
Set-CellText $ws 'D11' '0.07979'
Set-CellText $ws 'E11' '7.95%'
Set-CellText $ws 'D12' '0.08742'
Set-CellText $ws 'E12' '-2.07%'
Set-CellText $ws 'D13' '0.03130'
Set-CellText $ws 'E13' '3.08%'
Set-CellText $ws 'E14' '0.32%'
Set-CellText $ws 'D15' '0.001522'
Set-CellText $ws 'E15' '-0.43%'
Set-CellText $ws 'D16' '0.005732'
Set-CellText $ws 'E16' '-5.76%'
Set-CellText $ws 'B17' 'UpBots'
Set-CellText $ws 'C17' 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
Set-CellText $ws 'D17' '0.007500'
Set-CellText $ws 'E17' '2,097.97%'
Set-CellText $ws 'B18' 'LEO'
Set-CellText $ws 'C18' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-CellText $ws 'D18' '3.464'
Set-CellText $ws 'E18' '-3.19%'
Set-CellText $ws 'B19' 'BTSEToken'
Set-CellText $ws 'C19' 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-CellText $ws 'D19' '2.276'
Set-CellText $ws 'E19' '-0.44%'
Set-CellText $ws 'B20' 'BitpandaEcosystemToken'
Set-CellText $ws 'C20' 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-CellText $ws 'D20' '0.3288'
Set-CellText $ws 'E20' '-0.07%'
Set-CellText $ws 'B21' 'ProBitToken'
Set-CellText $ws 'C21' 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-CellText $ws 'D21' '0.1292'
Set-CellText $ws 'E21' '-3.54%'
Set-CellText $ws 'B22' 'MCDex'
Set-CellText $ws 'C22' 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-CellText $ws 'D22' '4.329'
Set-CellText $ws 'E22' '0.05%'
Set-CellText $ws 'B23' 'ZBToken'
Set-CellText $ws 'C23' 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
Set-CellText $ws 'D23' '0.1796'
Set-CellText $ws 'E23' '6.54%'
Set-CellText $ws 'B24' 'CoinExToken'
Set-CellText $ws 'C24' 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-CellText $ws 'D24' '0.04609'
Set-CellText $ws 'E24' '-0.38%'
Set-CellText $ws 'B25' 'BitKan'
Set-CellText $ws 'C25' 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-CellText $ws 'D25' '0.001238'
Set-CellText $ws 'E25' '-0.24%'
Set-CellText $ws 'B26' 'HotbitToken'
Set-CellText $ws 'C26' 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-CellText $ws 'D26' '0.004463'
Set-CellText $ws 'E26' '-1.55%'
Set-CellText $ws 'B27' 'NitroEx'
Set-CellText $ws 'C27' 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
Set-CellText $ws 'D27' '0.0001250'
Set-CellText $ws 'E27' '3.97%'
Set-CellText $ws 'B28' 'Spectre.aiUtilityToken'
Set-CellText $ws 'C28' 'https://coinranking.com/coin/yxQ8LoZvwJ6Ee+spectreaiutilitytoken-sxut'
Set-CellText $ws 'B29' 'LegolasExchange'
Set-CellText $ws 'C29' 'https://coinranking.com/coin/zEMEnlPs_94tc+legolasexchange-lgo'
Set-CellText $ws 'B30' 'BitZToken'
Set-CellText $ws 'C30' 'https://coinranking.com/coin/nLHHPBKDJSEee+bitztoken-bz'
Set-CellText $ws 'B31' 'Birake'
Set-CellText $ws 'C31' 'https://coinranking.com/coin/dTOfofFqKQiY5+birake-bir'
Set-CellText $ws 'B32' 'NashExchange'
Set-CellText $ws 'C32' 'https://coinranking.com/coin/9LcSTo0q-+nashexchange-nex'
Set-CellText $ws 'B33' 'AAXToken'
Set-CellText $ws 'C33' 'https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab'
Set-CellText $ws 'B34' 'CenX'
Set-CellText $ws 'C34' 'https://coinranking.com/coin/V4XJUvLQb+cenx-cenx'
Set-CellText $ws 'B35' 'BNIXToken'
Set-CellText $ws 'C35' 'https://coinranking.com/coin/n194X9uHp+bnixtoken-bnix'
Set-CellText $ws 'D39' '0.01724'
Set-CellText $ws 'E39' '-2.11%'
Set-CellText $ws 'D40' '0.04807'
Set-CellText $ws 'E40' '4.63%'
Set-CellText $ws 'D41' '0.007497'
Set-CellText $ws 'E41' '8.14%'
Set-CellText $ws 'E42' '-0.97%'
Set-CellText $ws 'D43' '0.002379'
Set-CellText $ws 'E43' '12.66%'
Set-CellText $ws 'D44' '0.01026'
Set-CellText $ws 'E44' '4.08%'
Set-CellText $ws 'D45' '0.00005974'
Set-CellText $ws 'E45' '-3.75%'
Set-CellText $ws 'D46' '0.00000000751'
Set-CellText $ws 'E46' '0.08%'
Set-CellText $ws 'D47' '0.003395'
Set-CellText $ws 'E47' '-59.62%'
Set-CellText $ws 'D48' '0.8204'
Set-CellText $ws 'E48' '2.36%'
Set-CellText $ws 'D49' '0.00002103'
Set-CellText $ws 'E49' '0.08%'
Set-CellText $ws 'D50' '0.0002003'
Set-CellText $ws 'E50' '0.08%'

$excel.CutCopyMode = 0
